$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Hoang Quoc Viet" profile (row 2: Name/Age/Gender/Plate/Phone) is
# deleted -- remaining rows shift up, matching the commit message
# "... Delete profile".
$ws.Rows.Item(2).Delete()

# Header for column F changes from "Table" to "ID".
$ws.Range("F1").Value = "ID"

# Font across the sheet moves from Arial to Calibri.
$ws.Range("A1:F1").Font.Name = "Calibri"
$ws.Range("A2:E3").Font.Name = "Calibri"
$ws.Range("A2:E3").ClearFormats()
